$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'Bel appartement neuf avec jardin'
$ws.Range("B2").Value = 'CHF 1,795,000.–'
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = '4'
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = '133m²'
$ws.Range("E2").Value = '1290 Versoix'
$ws.Range("F2").Value = 'https://www.homegate.ch/buy/4001266587'
$ws.Range("G2").Value = '2024-08-15 23:22:26'

# Row 3
$ws.Range("A3").Value = 'Appartement 5 pièces au 1er étage avec jardin'
$ws.Range("B3").Value = 'CHF 1,750,000.–'
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = '5'
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = '135m²'
$ws.Range("E3").Value = '1256 Troinex'
$ws.Range("F3").Value = 'https://www.homegate.ch/buy/4001276348'
$ws.Range("G3").Value = '2024-08-15 23:22:27'

# Row 4
$ws.Range("A4").Value = 'Chantier ouvert - Dernier lot disponible - Résidence Austra - 4 pièces en attique'
$ws.Range("B4").Value = 'CHF 1,095,000.–'
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = '4'
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = '74m²'
$ws.Range("E4").Value = 'Chemin Plein-sud 27, 1226 Thônex'
$ws.Range("F4").Value = 'https://www.homegate.ch/buy/4001358597'
$ws.Range("G4").Value = '2024-08-15 23:22:28'

# Row 5
$ws.Range("A5").Value = 'Dernier Lot - Livraison imminente - 4 pièces en attique - Éligible CASATAX'
$ws.Range("B5").Value = 'CHF 1,290,000.–'
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = '4'
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = '91m²'
$ws.Range("E5").Value = 'Route de Sous-Moulin 38A, 1226 Thônex'
$ws.Range("F5").Value = 'https://www.homegate.ch/buy/4000618526'
$ws.Range("G5").Value = '2024-08-15 23:22:29'

# Row 6
$ws.Range("A6").Value = 'Exclusivité! Magnifique attique neuf de 6 pièces avec terrasse d''angle'
$ws.Range("B6").Value = 'CHF 2,840,000.–'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = '6'
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = '192m²'
$ws.Range("E6").Value = '1256 Troinex'
$ws.Range("F6").Value = 'https://www.homegate.ch/buy/4001367678'
$ws.Range("G6").Value = '2024-08-15 23:22:30'

# Row 7
$ws.Range("A7").Value = 'Rare à la vente : Projet neuf disponible immédiatement'
$ws.Range("B7").Value = 'CHF 2,030,000.–'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = '5'
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = '139m²'
$ws.Range("E7").Value = '1228 Plan-les-Ouates'
$ws.Range("F7").Value = 'https://www.homegate.ch/buy/4001275870'
$ws.Range("G7").Value = '2024-08-15 23:22:31'

# Row 8
$ws.Range("A8").Value = 'LA LAC TOWER : nouvelle promotion immobilière de 18 appartements à vendre au Eaux-Vives'
$ws.Range("B8").Value = 'Price on request'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = '4.5'
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = 'N/A'
$ws.Range("E8").Value = 'Av. de Chamonix 3BIS, 1207 Genève'
$ws.Range("F8").Value = 'https://www.homegate.ch/buy/3001941475'
$ws.Range("G8").Value = '2024-08-15 23:22:32'

# Row 9
$ws.Range("A9").Value = 'Appartement avec un emplacement incroyable à Genève'
$ws.Range("B9").Value = 'CHF 1,550.–'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = '2.5'
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = '55m²'
$ws.Range("E9").Value = 'Rue de la Dôle 24, 1203 Genève'
$ws.Range("F9").Value = 'https://www.homegate.ch/buy/4001389730'
$ws.Range("G9").Value = '2024-08-15 23:22:33'

# Row 10
$ws.Range("A10").Value = 'Vieille Ville - Rue des Granges'
$ws.Range("B10").Value = 'CHF 4,575,000.–'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = '6'
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = '220m²'
$ws.Range("E10").Value = 'Rue des Granges, 1204 Genève'
$ws.Range("F10").Value = 'https://www.homegate.ch/buy/4001380333'
$ws.Range("G10").Value = '2024-08-15 23:22:35'

# Row 11
$ws.Range("A11").Value = 'Bel appartement traversant Champel'
$ws.Range("B11").Value = 'CHF 2,950,000.–'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = '6.5'
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = '179m²'
$ws.Range("E11").Value = 'Avenue de Champel 55, 1206 Genève'
$ws.Range("F11").Value = 'https://www.homegate.ch/buy/4001364526'
$ws.Range("G11").Value = '2024-08-15 23:22:36'

# Row 12
$ws.Range("A12").Value = 'Magnifique appartement neuf en rez-de-jardin au centre de Vésenaz'
$ws.Range("B12").Value = 'CHF 3,190,000.–'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = '5'
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = 'N/A'
$ws.Range("E12").Value = '6 Route de Vandoeuvres, 1222 Vésenaz'
$ws.Range("F12").Value = 'https://www.homegate.ch/buy/4001358835'
$ws.Range("G12").Value = '2024-08-15 23:22:37'

# Row 13
$ws.Range("A13").Value = 'Au coeur d''un quartier résidentiel!'
$ws.Range("B13").Value = 'Price on request'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = '7'
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = '200m²'
$ws.Range("E13").Value = '1223 Cologny'
$ws.Range("F13").Value = 'https://www.homegate.ch/buy/4001356168'
$ws.Range("G13").Value = '2024-08-15 23:22:38'

# Row 14
$ws.Range("A14").Value = 'Prenez de la hauteur!'
$ws.Range("B14").Value = 'Price on request'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = '8.5'
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = '200m²'
$ws.Range("E14").Value = '1203 Genève'
$ws.Range("F14").Value = 'https://www.homegate.ch/buy/4001356163'
$ws.Range("G14").Value = '2024-08-15 23:22:39'

# Row 15
$ws.Range("A15").Value = 'Exclusivité : appartement lumineux - dernier étage - quiétude et verdure'
$ws.Range("B15").Value = 'CHF 1,350,000.–'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = '5'
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = '118m²'
$ws.Range("E15").Value = 'Chemin du Daru, 1228 Plan-les-Ouates'
$ws.Range("F15").Value = 'https://www.homegate.ch/buy/4001317683'
$ws.Range("G15").Value = '2024-08-15 23:22:40'

# Row 16
$ws.Range("A16").Value = 'Appartement - Résidence située au bord du Lac - Lovée dans un vaste parc privatif de 1.4 hectares'
$ws.Range("B16").Value = 'CHF 1,790,000.–'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = '4'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = 'N/A'
$ws.Range("E16").Value = '1290 Versoix'
$ws.Range("F16").Value = 'https://www.homegate.ch/buy/4001315558'
$ws.Range("G16").Value = '2024-08-15 23:22:41'

# Row 17
$ws.Range("A17").Value = 'Studio à Vendre Champel/Florissant'
$ws.Range("B17").Value = 'CHF 640,000.–'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = '1.5'
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = '33m²'
$ws.Range("E17").Value = 'Chemin Rieu, 1208 Genève'
$ws.Range("F17").Value = 'https://www.homegate.ch/buy/4001270528'
$ws.Range("G17").Value = '2024-08-15 23:22:42'

# Row 18
$ws.Range("A18").Value = 'Bel appartement en duplex avec balcon et véranda'
$ws.Range("B18").Value = 'CHF 1,250,000.–'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = '5'
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = '130m²'
$ws.Range("E18").Value = 'Route d''Ambilly 36, 1226 Thônex'
$ws.Range("F18").Value = 'https://www.homegate.ch/buy/4001129344'
$ws.Range("G18").Value = '2024-08-15 23:22:43'

# Row 19
$ws.Range("A19").Value = 'Un magnifique duplex en attique'
$ws.Range("B19").Value = 'CHF 7,900,000.–'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = '12'
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = '521m²'
$ws.Range("E19").Value = 'Rue de la Cloche 6, 1201 Genève'
$ws.Range("F19").Value = 'https://www.homegate.ch/buy/3002187300'
$ws.Range("G19").Value = '2024-08-15 23:22:44'

# Row 20
$ws.Range("A20").Value = 'Tour Invictus : Appartement de 3.5 pièces à Fribourg avec vue panoramique sur les Préalpes !'
$ws.Range("B20").Value = 'CHF 799,000.–'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = '3.5'
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = '72m²'
$ws.Range("E20").Value = 'Route de la Glâne 19, 1202 Genève'
$ws.Range("F20").Value = 'https://www.homegate.ch/buy/4001319841'
$ws.Range("G20").Value = '2024-08-15 23:22:45'

# Row 21
$ws.Range("A21").Value = 'Entre le parc Bertrand et les HUG'
$ws.Range("B21").Value = 'CHF 1,670,000.–'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = '5.5'
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = '108m²'
$ws.Range("E21").Value = 'Genève, 1206 Genève'
$ws.Range("F21").Value = 'https://www.homegate.ch/buy/4001374783'
$ws.Range("G21").Value = '2024-08-15 23:22:46'
